# Weekly update: insert 4 new rows of "Melón - Tuna" price data (week of
# 2022-01-17) right after the current last "Tuna/Segunda" row (row 491),
# pushing the rest of the "Melón" block down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 492 - this shifts old rows 492:520 down to
# 496:524 (and copies the row-above formatting, matching Excel's default
# "insert" behaviour, e.g. the date-style on column D).
$ws.Rows("492:495").Insert()

# Row 492: Tuna / Extra
$ws.Cells.Item(492, 1).Value2 = 9
$ws.Cells.Item(492, 2).Value2 = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(492, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(492, 4).Value2 = 44578
$ws.Cells.Item(492, 5).Value2 = 13
$ws.Cells.Item(492, 6).Value2 = 100112027
$ws.Cells.Item(492, 7).Value2 = 'Melón'
$ws.Cells.Item(492, 8).Value2 = 'Tuna'
$ws.Cells.Item(492, 9).Value2 = 'Extra'
$ws.Cells.Item(492, 10).Value2 = 160
$ws.Cells.Item(492, 11).Value2 = 800
$ws.Cells.Item(492, 12).Value2 = 900
$ws.Cells.Item(492, 13).Value2 = 850
$ws.Cells.Item(492, 14).Value2 = '$/unidad'
$ws.Cells.Item(492, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(492, 16).Value2 = 850
$ws.Cells.Item(492, 17).Value2 = 1
$ws.Cells.Item(492, 18).Value2 = 'Hortaliza'

# Row 493: Tuna / Primera
$ws.Cells.Item(493, 1).Value2 = 9
$ws.Cells.Item(493, 2).Value2 = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(493, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(493, 4).Value2 = 44578
$ws.Cells.Item(493, 5).Value2 = 13
$ws.Cells.Item(493, 6).Value2 = 100112027
$ws.Cells.Item(493, 7).Value2 = 'Melón'
$ws.Cells.Item(493, 8).Value2 = 'Tuna'
$ws.Cells.Item(493, 9).Value2 = 'Primera'
$ws.Cells.Item(493, 10).Value2 = 250
$ws.Cells.Item(493, 11).Value2 = 600
$ws.Cells.Item(493, 12).Value2 = 700
$ws.Cells.Item(493, 13).Value2 = 650
$ws.Cells.Item(493, 14).Value2 = '$/unidad'
$ws.Cells.Item(493, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(493, 16).Value2 = 650
$ws.Cells.Item(493, 17).Value2 = 1
$ws.Cells.Item(493, 18).Value2 = 'Hortaliza'

# Row 494: Tuna / Segunda
$ws.Cells.Item(494, 1).Value2 = 9
$ws.Cells.Item(494, 2).Value2 = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(494, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(494, 4).Value2 = 44578
$ws.Cells.Item(494, 5).Value2 = 13
$ws.Cells.Item(494, 6).Value2 = 100112027
$ws.Cells.Item(494, 7).Value2 = 'Melón'
$ws.Cells.Item(494, 8).Value2 = 'Tuna'
$ws.Cells.Item(494, 9).Value2 = 'Segunda'
$ws.Cells.Item(494, 10).Value2 = 97
$ws.Cells.Item(494, 11).Value2 = 400
$ws.Cells.Item(494, 12).Value2 = 500
$ws.Cells.Item(494, 13).Value2 = 449
$ws.Cells.Item(494, 14).Value2 = '$/unidad'
$ws.Cells.Item(494, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(494, 16).Value2 = 449
$ws.Cells.Item(494, 17).Value2 = 1
$ws.Cells.Item(494, 18).Value2 = 'Hortaliza'

# Row 495: Tuna / Tercera
$ws.Cells.Item(495, 1).Value2 = 9
$ws.Cells.Item(495, 2).Value2 = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(495, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(495, 4).Value2 = 44578
$ws.Cells.Item(495, 5).Value2 = 13
$ws.Cells.Item(495, 6).Value2 = 100112027
$ws.Cells.Item(495, 7).Value2 = 'Melón'
$ws.Cells.Item(495, 8).Value2 = 'Tuna'
$ws.Cells.Item(495, 9).Value2 = 'Tercera'
$ws.Cells.Item(495, 10).Value2 = 52
$ws.Cells.Item(495, 11).Value2 = 300
$ws.Cells.Item(495, 12).Value2 = 300
$ws.Cells.Item(495, 13).Value2 = 300
$ws.Cells.Item(495, 14).Value2 = '$/unidad'
$ws.Cells.Item(495, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(495, 16).Value2 = 300
$ws.Cells.Item(495, 17).Value2 = 1
$ws.Cells.Item(495, 18).Value2 = 'Hortaliza'
